$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing header text in D1
$ws.Range("D1").Value = "CAGR RECEITAS 5 ANOS"

# Add new header in E1, matching style of D1 (copy formats, then set text)
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E1").Value = "CAGR LUCROS 5 ANOS"

# Fill E2:E21 with CAGR LUCROS 5 ANOS values
$values = @(
    "34,44%",
    "46,18%",
    "2,56%",
    "3,10%",
    "-%",
    "-26,63%",
    "-18,77%",
    "40,17%",
    "12,88%",
    "20,02%",
    "47,94%",
    "15,72%",
    "-%",
    "24,47%",
    "14,10%",
    "-%",
    "56,62%",
    "-%",
    "56,62%",
    "54,84%"
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 5).Value = $val
    $row++
}
